# Auto-generated update of market/profit data columns (H:N) for leve-crafting
# sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.
# Mirrors a scheduled data-refresh commit: numeric values in columns
# H (currentAveragePrice) .. N (LeveProfitHQ) are rewritten per-row;
# some rows gain/lose a trailing N (LeveProfitHQ) or M (LeveProfitNQ) cell
# when that quantity was/is not applicable for that row's leve.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1489.375
$ws.Range("I15").Value = 1489.375
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 4468.125
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -4299.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 3705470.2
$ws.Range("I88").Value = 5555705.5
$ws.Range("J88").Value = 5000
$ws.Range("K88").Value = 5555705.5
$ws.Range("L88").Value = 5000
$ws.Range("M88").Value = -5555299.5
$ws.Range("N88").Value = -5812

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 3705470.2
$ws.Range("I91").Value = 5555705.5
$ws.Range("J91").Value = 5000
$ws.Range("K91").Value = 5555705.5
$ws.Range("L91").Value = 5000
$ws.Range("M91").Value = -5554301.5
$ws.Range("N91").Value = -7808

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2977821
$ws.Range("I132").Value = 3677308.2
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 11031924.6
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -11029394.6
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 3690.875
$ws.Range("I135").Value = 4624.8696
$ws.Range("J135").Value = 1304
$ws.Range("K135").Value = 41623.8264
$ws.Range("L135").Value = 11736
$ws.Range("M135").Value = -39088.8264
$ws.Range("N135").Value = -16806

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 776.25
$ws.Range("I137").Value = 688.4545000000001
$ws.Range("J137").Value = 969.4
$ws.Range("K137").Value = 2065.3635
$ws.Range("L137").Value = 2908.2
$ws.Range("M137").Value = 484.6364999999996
$ws.Range("N137").Value = -8008.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2176.0562
$ws.Range("I138").Value = 721.41815
$ws.Range("J138").Value = 4529.147
$ws.Range("K138").Value = 2164.25445
$ws.Range("L138").Value = 13587.441
$ws.Range("M138").Value = 2975.74555
$ws.Range("N138").Value = -23867.441

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5480.7144
$ws.Range("I2").Value = 5480.7144
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 5480.7144
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -5367.7144
$ws.Range("N2").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6474.3
$ws.Range("I32").Value = 3793.2166
$ws.Range("J32").Value = 22560.8
$ws.Range("K32").Value = 3793.2166
$ws.Range("L32").Value = 22560.8
$ws.Range("M32").Value = -3506.2166
$ws.Range("N32").Value = -23134.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2600
$ws.Range("I61").Value = 2300
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 2300
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -2088
$ws.Range("N61").Value = -3424

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1029.25
$ws.Range("I74").Value = 1104.8572
$ws.Range("J74").Value = 500
$ws.Range("K74").Value = 1104.8572
$ws.Range("L74").Value = 500
$ws.Range("M74").Value = -230.8571999999999
$ws.Range("N74").Value = -2248

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1029.25
$ws.Range("I77").Value = 1104.8572
$ws.Range("J77").Value = 500
$ws.Range("K77").Value = 5524.286
$ws.Range("L77").Value = 2500
$ws.Range("M77").Value = -1156.286
$ws.Range("N77").Value = -11236

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2166.889
$ws.Range("I110").Value = 953.6667
$ws.Range("J110").Value = 4593.3335
$ws.Range("K110").Value = 953.6667
$ws.Range("L110").Value = 4593.3335
$ws.Range("M110").Value = 1091.3333
$ws.Range("N110").Value = -8683.333500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 5480.7144
$ws.Range("I116").Value = 5480.7144
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 5480.7144
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -3186.7144
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2119.625
$ws.Range("I122").Value = 1762.3636
$ws.Range("J122").Value = 2905.6
$ws.Range("K122").Value = 5287.0908
$ws.Range("L122").Value = 8716.799999999999
$ws.Range("M122").Value = -2837.0908
$ws.Range("N122").Value = -13616.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3438.7144
$ws.Range("I132").Value = 2086.182
$ws.Range("J132").Value = 4926.5
$ws.Range("K132").Value = 6258.545999999999
$ws.Range("L132").Value = 14779.5
$ws.Range("M132").Value = -3728.545999999999
$ws.Range("N132").Value = -19839.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2600
$ws.Range("I136").Value = 2300
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 6900
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -4350
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5480.7144
$ws.Range("I3").Value = 5480.7144
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5480.7144
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -5366.7144
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 539
$ws.Range("I80").Value = 104
$ws.Range("J80").Value = 653.4737
$ws.Range("K80").Value = 104
$ws.Range("L80").Value = 653.4737
$ws.Range("M80").Value = 894
$ws.Range("N80").Value = -2649.4737

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 539
$ws.Range("I83").Value = 104
$ws.Range("J83").Value = 653.4737
$ws.Range("K83").Value = 520
$ws.Range("L83").Value = 3267.3685
$ws.Range("M83").Value = 4472
$ws.Range("N83").Value = -13251.3685

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2878.56
$ws.Range("I134").Value = 2743.7144
$ws.Range("J134").Value = 3050.182
$ws.Range("K134").Value = 8231.143199999999
$ws.Range("L134").Value = 9150.545999999998
$ws.Range("M134").Value = -5696.143199999999
$ws.Range("N134").Value = -14220.546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 169466.67
$ws.Range("I3").Value = 250200
$ws.Range("J3").Value = 8000
$ws.Range("K3").Value = 250200
$ws.Range("L3").Value = 8000
$ws.Range("M3").Value = -250087
$ws.Range("N3").Value = -8226

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 9663.333000000001
$ws.Range("I16").Value = 7326.6665
$ws.Range("J16").Value = 12000
$ws.Range("K16").Value = 7326.6665
$ws.Range("L16").Value = 12000
$ws.Range("M16").Value = -7039.6665
$ws.Range("N16").Value = -12574

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1813.6216
$ws.Range("I31").Value = 1373.037
$ws.Range("J31").Value = 3003.2
$ws.Range("K31").Value = 1373.037
$ws.Range("L31").Value = 3003.2
$ws.Range("M31").Value = -1078.037

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1813.6216
$ws.Range("I34").Value = 1373.037
$ws.Range("J34").Value = 3003.2
$ws.Range("K34").Value = 1373.037
$ws.Range("L34").Value = 3003.2
$ws.Range("M34").Value = -1171.037

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2194.8235
$ws.Range("I99").Value = 2082
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 2082
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = -584
$ws.Range("N99").Value = -6996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 9663.333000000001
$ws.Range("I113").Value = 7326.6665
$ws.Range("J113").Value = 12000
$ws.Range("K113").Value = 7326.6665
$ws.Range("L113").Value = 12000
$ws.Range("M113").Value = -5156.6665
$ws.Range("N113").Value = -16340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 14227.3
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 14227.3
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 14227.3
$ws.Range("N125").Value = -19147.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2194.8235
$ws.Range("I126").Value = 2082
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 6246
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -3776
$ws.Range("N126").Value = -16940

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1680.6666
$ws.Range("I134").Value = 1515.3182
$ws.Range("J134").Value = 3499.5
$ws.Range("K134").Value = 4545.9546
$ws.Range("L134").Value = 10498.5
$ws.Range("M134").Value = -2010.9546
$ws.Range("N134").Value = -15568.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 73.117645
$ws.Range("I33").Value = 48.666668
$ws.Range("J33").Value = 100.625
$ws.Range("K33").Value = 292.000008
$ws.Range("L33").Value = 603.75
$ws.Range("M33").Value = -9.00000799999998
$ws.Range("N33").Value = -1169.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 39477200
$ws.Range("I121").Value = 500
$ws.Range("J121").Value = 40544140
$ws.Range("K121").Value = 1500
$ws.Range("L121").Value = 121632420
$ws.Range("M121").Value = -190
$ws.Range("N121").Value = -121635040

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1085.3125
$ws.Range("I129").Value = 300
$ws.Range("J129").Value = 1197.5
$ws.Range("K129").Value = 900
$ws.Range("L129").Value = 3592.5
$ws.Range("M129").Value = 4100
$ws.Range("N129").Value = -13592.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 758.5
$ws.Range("I131").Value = 393.52942
$ws.Range("J131").Value = 1172.1333
$ws.Range("K131").Value = 1180.58826
$ws.Range("L131").Value = 3516.3999
$ws.Range("M131").Value = 3859.41174
$ws.Range("N131").Value = -13596.3999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 8732.5
$ws.Range("I133").Value = 9465
$ws.Range("J133").Value = 8000
$ws.Range("K133").Value = 28395
$ws.Range("L133").Value = 24000
$ws.Range("M133").Value = -23335
$ws.Range("N133").Value = -34120

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 34675476
$ws.Range("I137").Value = 62513770
$ws.Range("J137").Value = 2860285.8
$ws.Range("K137").Value = 187541310
$ws.Range("L137").Value = 8580857.399999999
$ws.Range("M137").Value = -187536210
$ws.Range("N137").Value = -8591057.399999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2802.125
$ws.Range("I122").Value = 2283.4
$ws.Range("J122").Value = 3666.6667
$ws.Range("K122").Value = 6850.200000000001
$ws.Range("L122").Value = 11000.0001
$ws.Range("M122").Value = -4400.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5285.355
$ws.Range("I132").Value = 9145.286
$ws.Range("J132").Value = 2106.5881
$ws.Range("K132").Value = 27435.858
$ws.Range("L132").Value = 6319.7643
$ws.Range("M132").Value = -24905.858
$ws.Range("N132").Value = -11379.7643

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2274.4524
$ws.Range("I132").Value = 1478.75
$ws.Range("J132").Value = 3865.8572
$ws.Range("K132").Value = 4436.25
$ws.Range("L132").Value = 11597.5716
$ws.Range("M132").Value = -1906.25
$ws.Range("N132").Value = -16657.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3666.3333
$ws.Range("I122").Value = 2499.5
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 7498.5
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -5048.5
$ws.Range("N122").Value = -22900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2539.923
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2539.923
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 7619.768999999999
$ws.Range("N132").Value = -12679.769
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4236.8486
$ws.Range("I136").Value = 1084.0526
$ws.Range("J136").Value = 8515.643
$ws.Range("K136").Value = 3252.1578
$ws.Range("L136").Value = 25546.929
$ws.Range("M136").Value = -702.1578
$ws.Range("N136").Value = -30646.929
